# B6-PowerPoint.pptx — theme swap (Integral -> Office Theme) + table style update
# Commit: Mon, Jun 01, 2020  9:05:27 AM

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Re-colour the shared theme (ppt/theme/theme1.xml) from the "Integral" /
#    "Red Violet" palette to the default "Office Theme" palette. All slides
#    share the single slide master/theme in this deck, so editing the
#    master's legacy 8/12-slot ColorScheme rewrites theme1.xml directly.
#    Index -> theme role: 1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2
#    7=accent3 8=accent4 9=accent5 10=accent6 11=hlink 12=folHlink
#    COM RGB integers are packed 0x00BBGGRR.
# ---------------------------------------------------------------------------
$master = $p.Slides.Item(1).Master
$scheme = $master.ColorScheme

$scheme.Colors(1).RGB  = 0        # dk1      000000 (unchanged)
$scheme.Colors(2).RGB  = 16777215 # lt1      FFFFFF (unchanged)
$scheme.Colors(3).RGB  = 6968388  # dk2      44546A
$scheme.Colors(4).RGB  = 15132391 # lt2      E7E6E6
$scheme.Colors(5).RGB  = 13998939 # accent1  5B9BD5
$scheme.Colors(6).RGB  = 3243501  # accent2  ED7D31
$scheme.Colors(7).RGB  = 10855845 # accent3  A5A5A5
$scheme.Colors(8).RGB  = 49407    # accent4  FFC000
$scheme.Colors(9).RGB  = 12874308 # accent5  4472C4
$scheme.Colors(10).RGB = 4697456  # accent6  70AD47
$scheme.Colors(11).RGB = 12673797 # hlink    0563C1
$scheme.Colors(12).RGB = 7491477  # folHlink 954F72

# ---------------------------------------------------------------------------
# 2) Point the three existing tables at the new built-in table style.
# ---------------------------------------------------------------------------
$newStyleId = "{C22BAA12-D315-4E26-A71B-8B5C797467DF}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    $table = $slide.Shapes.Item(1).Table
    $table.ApplyStyle($newStyleId)
}
